$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Day 14 (row 15) entered first
$ws.Range("B15").Value = "LS vor dem Schlafen"
$ws.Range("C15").Value = "In den letzten Wochen waren meine Nächte vom Sonntag auf den Montag nicht sehr gut - konnte nicht gut einschlafen. Diesmal habe ich vor dem Schlafen LS gemacht mit den Eindrücken vom Tag und mit den nervösen Gefühlen, die ich hatte. Ich konnte nachher rasch einschlafen."
$ws.Range("D15").Value = "LS anwenden"
$ws.Range("E15").Value = "https://greator.com/wp-content/uploads/2021/05/selbstcoaching-step-by-step-as-206759352-1024x683.jpeg"

# Day 13 (row 14) entered second
$ws.Range("B14").Value = "Genussvoller, entspannter Tag"
$ws.Range("C14").Value = "Am Samstag haben wir einen Geburtstag gefeiert mit einem Brunch und waren danach noch Abendessen mit Freunden. Es war ein entspannter Tag, den ich richtig geniessen konnte. Es gab 2, 3 Mal Momente, in denen ich wieder etwas weniger entspannt war. Ich konnte mich aber gut abgrenzen und die Energie, die nicht zu mir gehört, habe ich auch nicht angenommen."
$ws.Range("E14").Value = "https://static.wikia.nocookie.net/kardashev/images/0/0f/Energy.jpg/revision/latest/scale-to-width-down/1000?cb=20210905164631"
$ws.Range("D14").Value = "Abgrenzung"

$ws.Range("C19").Select()
